# Apply the "Updated cryptos list" data refresh to the crypto price table.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Row 1 is the header.
#
# Several new Price values are plain decimal numbers (e.g. "0.999", "6.10").
# Excel/IronCalc would normally auto-infer those as numeric cells, but the
# source workbook keeps every Price/Volume cell as literal text (inlineStr,
# no number format). To reproduce that faithfully we temporarily force the
# cell to Text format ("@") before writing the value, then restore the
# cell style back to "Normal" so no visible formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '64.238.17'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '3.081.35'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue "D5" '557.99'
$ws.Range("E5").Value = '  +1.17%  '
Set-TextValue "D6" '144.03'
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.080.79'
$ws.Range("E8").Value = '  +0.50%  '
Set-TextValue "D9" '0.508'
$ws.Range("E9").Value = '  +1.01%  '
Set-TextValue "D10" '0.156'
$ws.Range("E10").Value = '  +2.29%  '
Set-TextValue "D11" '6.10'
$ws.Range("E11").Value = '  -6.23%  '
$ws.Range("E12").Value = '  +3.25%  '
Set-TextValue "D13" '0.0000230'
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").Value = '3.595.51'
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").Value = '64.253.22'
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("D17").Value = '3.084.10'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("E18").Value = '  +1.16%  '
Set-TextValue "D19" '6.73'
$ws.Range("E19").Value = '  -0.87%  '
Set-TextValue "D20" '480.58'
$ws.Range("E20").Value = '  -0.96%  '
Set-TextValue "D21" '14.03'
$ws.Range("E21").Value = '  +0.84%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  +3.09%  '
Set-TextValue "D24" '14.11'
$ws.Range("E24").Value = '  +10.35%  '
Set-TextValue "D25" '81.42'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  +0.00%  '
Set-TextValue "D27" '2.79'
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("E28").Value = '  +1.09%  '
$ws.Range("E29").Value = '  +1.86%  '
Set-TextValue "D30" '0.999'
$ws.Range("E30").Value = '  -0.17%  '
Set-TextValue "D31" '26.30'
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("E33").Value = '  +0.68%  '
Set-TextValue "D34" '5.63'
$ws.Range("E34").Value = '  -1.38%  '
Set-TextValue "D35" '6.21'
$ws.Range("E35").Value = '  +3.20%  '
Set-TextValue "D36" '55.46'
$ws.Range("E36").Value = '  -0.17%  '
Set-TextValue "D37" '0.0408'
$ws.Range("E37").Value = '  +2.12%  '
$ws.Range("E38").Value = '  +14.37%  '
Set-TextValue "D39" '439.25'
$ws.Range("E39").Value = '  -6.09%  '
Set-TextValue "D40" '0.0809'
$ws.Range("E40").Value = '  -1.95%  '
$ws.Range("D41").Value = '2.959.31'
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("E43").Value = '  -4.53%  '
Set-TextValue "D44" '28.20'
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("E47").Value = '  +3.88%  '
$ws.Range("E48").Value = '  +1.21%  '
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.0₃0516'
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D50" '117.79'
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("E51").Value = '  -0.52%  '
